$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) ---
$ws.Range("A1").Value = "TABLE NAME"
$ws.Range("B1").Value = "ATTRUBITE NAME"
$ws.Range("C1").Value = "TYPE"

# --- Skin Issues table (rows 2-7) ---
$ws.Range("A2").Value = "Skin Issues"
$ws.Range("B2").Value = "Id"
$ws.Range("C2").Value = "INTEGER"

$ws.Range("B3").Value = "IssueName"
$ws.Range("C3").Value = "STRING"

$ws.Range("B4").Value = "IssueColor"
$ws.Range("C4").Value = "STRING"

$ws.Range("B5").Value = " IssueDefination"
$ws.Range("C5").Value = "STRING"

$ws.Range("B6").Value = "DisplayOrder"
$ws.Range("C6").Value = "INTEGER"

$ws.Range("B7").Value = "IsActive"
$ws.Range("C7").Value = "BOOLEAN"

# --- Score Categories table (rows 9-14) ---
$ws.Range("A9").Value = "Score Categories"
$ws.Range("B9").Value = "Id"
$ws.Range("C9").Value = "INTEGER"

$ws.Range("B10").Value = "CategoryName"
$ws.Range("C10").Value = "STRING"

$ws.Range("B11").Value = "CategoryColor"
$ws.Range("C11").Value = "STRING"

$ws.Range("B12").Value = "MinScore"
$ws.Range("C12").Value = "INTEGER"

$ws.Range("B13").Value = "MaxScore"
$ws.Range("C13").Value = "INTEGER"

$ws.Range("B14").Value = "IsActive"
$ws.Range("C14").Value = "BOOLEAN"

# --- User Info table (rows 16-18) ---
$ws.Range("A16").Value = "User Info"
$ws.Range("B16").Value = "Id"
$ws.Range("C16").Value = "INTEGER"

$ws.Range("B17").Value = "UserUuid"
$ws.Range("C17").Value = "STRING"

$ws.Range("B18").Value = "Source"
$ws.Range("C18").Value = "STRING"

# --- Score History table (rows 20-22) ---
$ws.Range("A20").Value = "Score History"
$ws.Range("B20").Value = "Id"
$ws.Range("C20").Value = "INTEGER"

$ws.Range("B21").Value = "UserUuid"
$ws.Range("C21").Value = "STRING"

$ws.Range("B22").Value = "Concerns"
$ws.Range("C22").Value = "JSON"

# --- blank separator rows ---
$ws.Range("A8").HorizontalAlignment = -4142
$ws.Range("A19").HorizontalAlignment = -4142

# --- merges ---
$ws.Range("A16:A18").Merge()

# --- column C alignment (center horizontally, like columns A/B) ---
$ws.Range("C1:C7").HorizontalAlignment = -4108
$ws.Range("C1:C7").VerticalAlignment = -4108
$ws.Range("C9:C14").HorizontalAlignment = -4108
$ws.Range("C9:C14").VerticalAlignment = -4108
$ws.Range("C16:C18").HorizontalAlignment = -4108
$ws.Range("C16:C18").VerticalAlignment = -4108
$ws.Range("C20:C22").HorizontalAlignment = -4108
$ws.Range("C20:C22").VerticalAlignment = -4108

# --- column widths ---
$ws.Columns.Item(1).ColumnWidth = 16
$ws.Columns.Item(2).ColumnWidth = 17.1818181818182
$ws.Columns.Item(3).ColumnWidth = 9.63636363636364

# --- sheet view ---
$ws.Range("F12").Select()
$excel.ActiveWindow.ScrollRow = 7

